$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3124.6365
$ws.Range("I28").Value = 2624.5557
$ws.Range("K28").Value = 2624.5557
$ws.Range("M28").Value = -2139.5557
$ws.Range("H137").Value = 23138.482
$ws.Range("I137").Value = 18554.117
$ws.Range("J137").Value = 30931.9
$ws.Range("K137").Value = 55662.351
$ws.Range("L137").Value = 92795.70000000001
$ws.Range("M137").Value = -53112.351
$ws.Range("N137").Value = -97895.70000000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3705.2173
$ws.Range("J63").Value = 4633.1333
$ws.Range("L63").Value = 4633.1333
$ws.Range("N63").Value = -6005.1333
$ws.Range("H66").Value = 3705.2173
$ws.Range("J66").Value = 4633.1333
$ws.Range("L66").Value = 23165.6665
$ws.Range("N66").Value = -30029.6665
$ws.Range("H74").Value = 679534.5600000001
$ws.Range("I74").Value = 1200799.6
$ws.Range("K74").Value = 1200799.6
$ws.Range("M74").Value = -1199925.6
$ws.Range("H77").Value = 679534.5600000001
$ws.Range("I77").Value = 1200799.6
$ws.Range("K77").Value = 6003998
$ws.Range("M77").Value = -5999630
$ws.Range("H97").Value = 1345.8438
$ws.Range("I97").Value = 1067.52
$ws.Range("J97").Value = 2339.8572
$ws.Range("K97").Value = 1067.52
$ws.Range("L97").Value = 2339.8572
$ws.Range("M97").Value = -571.52
$ws.Range("N97").Value = -3331.8572
$ws.Range("H132").Value = 1818.1538
$ws.Range("I132").Value = 1142
$ws.Range("K132").Value = 3426
$ws.Range("M132").Value = -896

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1880.5714
$ws.Range("I20").Value = 1805.2307
$ws.Range("J20").Value = 2003
$ws.Range("K20").Value = 1805.2307
$ws.Range("L20").Value = 2003
$ws.Range("M20").Value = -1558.2307
$ws.Range("N20").Value = -2497
$ws.Range("H64").Value = 1271.3334
$ws.Range("I64").Value = 1742.6666
$ws.Range("J64").Value = 1114.2222
$ws.Range("K64").Value = 1742.6666
$ws.Range("L64").Value = 1114.2222
$ws.Range("M64").Value = -1517.6666
$ws.Range("N64").Value = -1564.2222
$ws.Range("H67").Value = 1271.3334
$ws.Range("I67").Value = 1742.6666
$ws.Range("J67").Value = 1114.2222
$ws.Range("K67").Value = 1742.6666
$ws.Range("L67").Value = 1114.2222
$ws.Range("M67").Value = -962.6666
$ws.Range("N67").Value = -2674.2222
$ws.Range("H80").Value = 791.1875
$ws.Range("I80").Value = 924.75
$ws.Range("J80").Value = 657.625
$ws.Range("K80").Value = 924.75
$ws.Range("L80").Value = 657.625
$ws.Range("M80").Value = 73.25
$ws.Range("N80").Value = -2653.625
$ws.Range("H83").Value = 791.1875
$ws.Range("I83").Value = 924.75
$ws.Range("J83").Value = 657.625
$ws.Range("K83").Value = 4623.75
$ws.Range("L83").Value = 3288.125
$ws.Range("M83").Value = 368.25
$ws.Range("N83").Value = -13272.125

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5001350.5
$ws.Range("I31").Value = 5556445.5
$ws.Range("J31").Value = 5499.5
$ws.Range("K31").Value = 5556445.5
$ws.Range("L31").Value = 5499.5
$ws.Range("M31").Value = -5556150.5
$ws.Range("N31").Value = -6089.5
$ws.Range("H34").Value = 5001350.5
$ws.Range("I34").Value = 5556445.5
$ws.Range("J34").Value = 5499.5
$ws.Range("K34").Value = 5556445.5
$ws.Range("L34").Value = 5499.5
$ws.Range("M34").Value = -5556243.5
$ws.Range("N34").Value = -5903.5
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("H134").Value = 3095.0435
$ws.Range("I134").Value = 2672.9473
$ws.Range("K134").Value = 8018.841899999999
$ws.Range("M134").Value = -5483.841899999999
$ws.Range("N118").Value = -53314

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 108485
$ws.Range("I131").Value = 357719.1
$ws.Range("K131").Value = 1073157.3
$ws.Range("M131").Value = -1068117.3

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4978.778
$ws.Range("I70").Value = 4821.5
$ws.Range("J70").Value = 5293.3335
$ws.Range("K70").Value = 4821.5
$ws.Range("L70").Value = 5293.3335
$ws.Range("M70").Value = -4551.5
$ws.Range("N70").Value = -5833.3335
$ws.Range("H73").Value = 4978.778
$ws.Range("I73").Value = 4821.5
$ws.Range("J73").Value = 5293.3335
$ws.Range("K73").Value = 4821.5
$ws.Range("L73").Value = 5293.3335
$ws.Range("M73").Value = -3885.5
$ws.Range("N73").Value = -7165.3335
$ws.Range("H113").Value = 2997.3635
$ws.Range("I113").Value = 2796.647
$ws.Range("J113").Value = 3679.8
$ws.Range("K113").Value = 2796.647
$ws.Range("L113").Value = 3679.8
$ws.Range("M113").Value = -626.6469999999999
$ws.Range("N113").Value = -8019.8

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3404
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H136").Value = 4516.3076
$ws.Range("I136").Value = 4121
$ws.Range("K136").Value = 12363
$ws.Range("M136").Value = -9813
$ws.Range("M61").ClearContents()
$ws.Range("M113").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3926.25
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3926.25
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3926.25
$ws.Range("N62").Value = -5174.25
$ws.Range("H65").Value = 3926.25
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3926.25
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 19631.25
$ws.Range("N65").Value = -25871.25
$ws.Range("H136").Value = 19647.445
$ws.Range("I136").Value = 23424.482
$ws.Range("J136").Value = 3999.7144
$ws.Range("K136").Value = 70273.446
$ws.Range("L136").Value = 11999.1432
$ws.Range("M136").Value = -67723.446
$ws.Range("N136").Value = -17099.1432
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()
